# Add a "canonical SMILES" column (D) to the microstate list sheet.
# For this molecule the canonical SMILES is identical to the canonical
# isomeric SMILES already present in column C, so we duplicate those
# values into the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column D.
$ws.Cells.Item(2, 4).Value = "canonical SMILES"

# Duplicate the canonical isomeric SMILES (column C) into the new
# canonical SMILES column (D) for every data row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 3; $r -le $lastRow; $r++) {
    $smiles = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 4).Value = $smiles
}

# Give the new column a sensible width, matching the others (closest
# achievable value to 36.85546875 character-widths given COM rounding).
$ws.Columns.Item(4).ColumnWidth = 36
